$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so values that look like
# numbers (e.g. "1.002", "17.00") are stored as text, matching the source data.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2: update Price, Volume
$ws.Range("D2").Value = "28.913.27"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3: update Price, Volume
$ws.Range("D3").Value = "1.879.57"
$ws.Range("E3").Value = "  +0.00%  "

# Row 4: update Volume
$ws.Range("E4").Value = "  -0.72%  "

# Row 5: update Price, Volume
$ws.Range("D5").Value = "324.95"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6: update Price, Volume
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.63%  "

# Row 7: update Price, Volume
$ws.Range("D7").Value = "0.4604"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8: update Price, Volume
$ws.Range("D8").Value = "0.3877"
$ws.Range("E8").Value = "  +0.44%  "

# Row 9: update Price, Volume
$ws.Range("D9").Value = "0.07859"
$ws.Range("E9").Value = "  -0.01%  "

# Row 10: update Price, Volume
$ws.Range("D10").Value = "0.9863"
$ws.Range("E10").Value = "  -1.49%  "

# Row 11: update Price, Volume
$ws.Range("D11").Value = "21.77"
$ws.Range("E11").Value = "  +0.32%  "

# Row 12: update Coin, Link, Price, Volume
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.908.77"
$ws.Range("E12").Value = "  +1.64%  "

# Row 13: update Coin, Link, Price, Volume
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "6.999"
$ws.Range("E13").Value = "  -0.88%  "

# Row 14: update Coin, Link, Price, Volume
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.649"
$ws.Range("E14").Value = "  -0.98%  "

# Row 15: update Price, Volume
$ws.Range("D15").Value = "0.06953"
$ws.Range("E15").Value = "  -0.52%  "

# Row 16: update Price
$ws.Range("D16").Value = "88.06"

# Row 17: update Volume
$ws.Range("E17").Value = "  -0.60%  "

# Row 18: update Volume
$ws.Range("E18").Value = "  -0.56%  "

# Row 19: update Price, Volume
$ws.Range("D19").Value = "17.00"
$ws.Range("E19").Value = "  -1.10%  "

# Row 20: update Price, Volume
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.71%  "

# Row 21: update Price, Volume
$ws.Range("D21").Value = "28.891.40"
$ws.Range("E21").Value = "  +1.11%  "

# Row 22: update Price, Volume
$ws.Range("D22").Value = "5.233"
$ws.Range("E22").Value = "  -1.87%  "

# Row 23: update Price, Volume
$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  -0.18%  "

# Row 24: update Price, Volume
$ws.Range("D24").Value = "2.086"
$ws.Range("E24").Value = "  +1.34%  "

# Row 25: update Price, Volume
$ws.Range("D25").Value = "156.34"
$ws.Range("E25").Value = "  +1.26%  "

# Row 26: update Price, Volume
$ws.Range("D26").Value = "19.30"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27: update Price, Volume
$ws.Range("D27").Value = "5.999"
$ws.Range("E27").Value = "  +2.53%  "

# Row 28: update Price, Volume
$ws.Range("D28").Value = "1.927"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29: update Price, Volume
$ws.Range("D29").Value = "117.46"
$ws.Range("E29").Value = "  -0.77%  "

# Row 30: update Price, Volume
$ws.Range("D30").Value = "0.09368"
$ws.Range("E30").Value = "  +0.26%  "

# Row 31: update Volume
$ws.Range("E31").Value = "  -1.98%  "

# Row 32: update Price, Volume
$ws.Range("D32").Value = "5.261"
$ws.Range("E32").Value = "  -0.96%  "

# Row 33: update Price, Volume
$ws.Range("D33").Value = "1.313"
$ws.Range("E33").Value = "  -1.64%  "

# Row 34: update Volume
$ws.Range("E34").Value = "  -0.55%  "

# Row 35: update Volume
$ws.Range("E35").Value = "  +1.22%  "

# Row 36: update Volume
$ws.Range("E36").Value = "  -0.58%  "

# Row 37: update Price, Volume
$ws.Range("D37").Value = "0.02071"
$ws.Range("E37").Value = "  +0.12%  "

# Row 38: update Volume
$ws.Range("E38").Value = "  -0.71%  "

# Row 39: update Price, Volume
$ws.Range("D39").Value = "7.618"
$ws.Range("E39").Value = "  -4.59%  "

# Row 40: update Price, Volume
$ws.Range("D40").Value = "0.5648"
$ws.Range("E40").Value = "  -0.54%  "

# Row 41: update Volume
$ws.Range("E41").Value = "  -1.61%  "

# Row 42: update Price, Volume
$ws.Range("D42").Value = "9.675"
$ws.Range("E42").Value = "  -0.57%  "

# Row 43: update Price, Volume
$ws.Range("D43").Value = "2.272"
$ws.Range("E43").Value = "  +5.80%  "

# Row 44: update Price, Volume
$ws.Range("D44").Value = "11.91"
$ws.Range("E44").Value = "  +1.26%  "

# Row 45: update Price
$ws.Range("D45").Value = "0.5345"

# Row 46: update Price
$ws.Range("D46").Value = "0.07045"

# Row 47: update Price, Volume
$ws.Range("D47").Value = "1.843"
$ws.Range("E47").Value = "  +0.15%  "

# Row 48: update Price, Volume
$ws.Range("D48").Value = "112.70"
$ws.Range("E48").Value = "  +0.35%  "

# Row 49: update Price, Volume
$ws.Range("D49").Value = "2.531"
$ws.Range("E49").Value = "  +1.48%  "

# Row 50: update Price, Volume
$ws.Range("D50").Value = "1.068"
$ws.Range("E50").Value = "  -4.57%  "

# Row 51: update Coin, Link, Price, Volume
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "70.85"
$ws.Range("E51").Value = "  +0.19%  "

# Restore column D to the default (unstyled) cell style now that the
# values are safely stored as text, so no stray number-format styling
# is left behind on the cells themselves.
$dRange.Style = "Normal"
